$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the URL (row 2) - matchsource -> matchsync
$ws.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-preferred-product-codes"

# Set the Experimental value (row 7) to the literal text "true" (not a Boolean).
# A plain .Value assignment of "true" is auto-coerced to a Boolean by Excel, so
# instead write it as a formula that evaluates to the text string, then convert
# that formula to its static value via copy / paste-special values, which keeps
# the cell's text type and original style.
$cB7 = $ws.Range("B7")
$cB7.Formula = "=""true"""
$cB7.Copy()
$cB7.PasteSpecial(-4163)

# Update the Date (row 8) to the new timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
